# [Fonds de solidarite] Add 2020-12-30 data
#
# The source data refresh bumps the "nombre_aides" (column C) and
# "montant_total" (column D) counters for a handful of (region x section)
# rows. Both columns are stored as plain text in the sheet (e.g. "130281.00"
# keeps its trailing zeroes), so every write below forces the cell to text
# format first and then restores the "Normal" style so no stray formatting
# is left behind - only the literal displayed text changes, exactly like in
# the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new nombre_aides (C), new montant_total (D)
$updates = @(
    @{ Row = 21;  C = "48";   D = "130281.00" },
    @{ Row = 24;  C = "494";  D = "3806709.26" },
    @{ Row = 39;  C = "578";  D = "6190211.97" },
    @{ Row = 47;  C = "96";   D = "1028987.72" },
    @{ Row = 54;  C = "393";  D = "3338101.20" },
    @{ Row = 58;  C = "51";   D = "422331.00" },
    @{ Row = 59;  C = "26";   D = "153262.00" },
    @{ Row = 84;  C = "830";  D = "7220201.36" },
    @{ Row = 130; C = "1125"; D = "9105363.67" },
    @{ Row = 155; C = "834";  D = "4341626.43" },
    @{ Row = 174; C = "95";   D = "359128.34" },
    @{ Row = 214; C = "1006"; D = "9724988.80" },
    @{ Row = 222; C = "208";  D = "2189031.09" },
    @{ Row = 245; C = "475";  D = "3465751.80" },
    @{ Row = 249; C = "104";  D = "809165.53" },
    @{ Row = 253; C = "92";   D = "818127.57" },
    @{ Row = 261; C = "1710"; D = "11976409.68" },
    @{ Row = 265; C = "293";  D = "1585454.81" },
    @{ Row = 269; C = "288";  D = "2057321.36" },
    @{ Row = 270; C = "323";  D = "837082.41" }
)

foreach ($u in $updates) {
    $cRng = $ws.Range("C" + $u.Row)
    $cRng.NumberFormat = "@"
    $cRng.Value = $u.C
    $cRng.Style = "Normal"

    $dRng = $ws.Range("D" + $u.Row)
    $dRng.NumberFormat = "@"
    $dRng.Value = $u.D
    $dRng.Style = "Normal"
}
